# Aggiorno file need_to_buy.xlsx da R
# Daily refresh of the forecast/buy table: the rolling date window shifts
# forward by one day and the forecast-derived columns (B,C,E,F) are
# recomputed; D is mostly unchanged except the first row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  2  = @(46030, 14502.5252860584, 14220.9239864974, 23835.86, 9624.96589432651, 0.417911700996228)
  3  = @(46031, 13772.887967387,  12280.67667569,   12075.86, 9471.91632533489, 403.197208376036)
  4  = @(46032, 5623.08451971876, 8779.43459101948, 12075.86, 8900.52940047133, 233.504332978783)
  5  = @(46033, 5501.22637197257, 9028.23850089512, 12075.86, 8971.5176833041,  246.829007674967)
  6  = @(46034, 13455.7123880986, 13405.2258534184, 12075.86, 9197.02100925178, 438.599452611257)
  7  = @(46035, 13378.423123507,  12846.0900743681, 12075.86, 9107.15683099638, 411.557787723522)
  8  = @(46036, 13378.423123507,  12416.1466958139, 12075.86, 9107.15683099638, 393.643480283761)
  9  = @(46037, 13378.423123507,  12333.0655568506, 12075.86, 9107.15683099638, 390.181766160291)
  10 = @(46038, 13378.423123507,  11495.1682058025, 12075.86, 9107.1167504224,  355.267706509371)
  11 = @(46039, 5471.00037786234, 7743.20935132198, 12075.86, 8692.30583471229, 181.652299418095)
  12 = @(46040, 5302.34398887746, 7577.57805082373, 12075.86, 8679.88993889267, 174.233666238183)
  13 = @(46041, 13071.6216038713, 11472.3479477444, 12075.86, 8882.63087916091, 344.963284454387)
  14 = @(46042, 13071.6216038713, 11671.5194650052, 12075.86, 8882.63087916091, 353.262097673589)
  15 = @(46043, 13071.6216038713, 11846.8933320964, 12075.86, 8882.63087916091, 360.569342135719)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $values[$i]
    }
}
